# "Generate List of Order" - populate the Size column for the Musical
# chairs mechanism on the "Morning Order" sheet, and update the last
# selected cell / window layout to reflect the saved state.

$wb = $excel.ActiveWorkbook

# --- Best-effort: restore the workbook window layout (maximized) ---
$win = $wb.Windows.Item(1)
$win.Left = -120
$win.Top = -120
$win.Width = 29040
$win.Height = 17640

# --- Work on the "Morning Order" sheet ---
$ws = $wb.Worksheets.Item("Morning Order")
$ws.Activate()

# Fill in the Size values (column C) for faculties 8-12 (rows 9-13),
# previously placeholder zeros.
$ws.Range("C9").Value = 303
$ws.Range("C10").Value = 103
$ws.Range("C11").Value = 411
$ws.Range("C12").Value = 221
$ws.Range("C13").Value = 305

# Leave the selection where the author left it when saving (E16).
$ws.Range("E16").Select()
